$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "57÷4=14, 1"
$t.Cell(1,2).Range.Text = "26÷8=3, 2"
$t.Cell(1,3).Range.Text = "55÷3=18, 1"
$t.Cell(1,4).Range.Text = "91÷5=18, 1"
$t.Cell(1,5).Range.Text = "36÷4=9, 0"
$t.Cell(5,1).Range.Text = "57÷6=9, 3"
$t.Cell(5,2).Range.Text = "51÷8=6, 3"
$t.Cell(5,3).Range.Text = "55÷8=6, 7"
$t.Cell(5,4).Range.Text = "94÷7=13, 3"
$t.Cell(5,5).Range.Text = "95÷7=13, 4"
$t.Cell(9,1).Range.Text = "90÷2=45, 0"
$t.Cell(9,2).Range.Text = "70÷9=7, 7"
$t.Cell(9,3).Range.Text = "67÷3=22, 1"
$t.Cell(9,4).Range.Text = "22÷2=11, 0"
$t.Cell(9,5).Range.Text = "21÷9=2, 3"
$t.Cell(13,1).Range.Text = "80÷2=40, 0"
$t.Cell(13,2).Range.Text = "52÷7=7, 3"
$t.Cell(13,3).Range.Text = "60÷4=15, 0"
$t.Cell(13,4).Range.Text = "13÷8=1, 5"
$t.Cell(13,5).Range.Text = "37÷5=7, 2"
$t.Cell(17,1).Range.Text = "57÷8=7, 1"
$t.Cell(17,2).Range.Text = "77÷5=15, 2"
$t.Cell(17,3).Range.Text = "25÷4=6, 1"
$t.Cell(17,4).Range.Text = "27÷8=3, 3"
$t.Cell(17,5).Range.Text = "71÷8=8, 7"
